$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simulated energy-consumption values (rows 2-49, columns B and C)
# using the faster (priority-queue-based) simulation results.
$values = @{
    2 = @(0.9136590783126551, 1.683700707322582)
    3 = @(3.602355127675703, 3.799928533123985)
    4 = @(4.583538481480849, 5.575007633379273)
    5 = @(6.116689810046933, 7.321119304305469)
    6 = @(6.874668290849217, 9.124756499179981)
    7 = @(12.64272738760497, 10.71236735005599)
    8 = @(12.73712352956607, 12.45353484821274)
    9 = @(15.1170697183871, 14.69169274034213)
    10 = @(17.30975015175841, 16.52036898181901)
    11 = @(19.82372282953638, 18.20514238734986)
    12 = @(20.10676859513826, 20.10837165973416)
    13 = @(21.53475903542832, 21.65304937383792)
    14 = @(22.84015723675692, 23.33378020994451)
    15 = @(28.45862209557449, 25.0982114602135)
    16 = @(32.78682863864541, 26.87610102991496)
    17 = @(33.63334990107052, 28.83095017427892)
    18 = @(36.31181953000279, 30.51678714007469)
    19 = @(37.3120120359229, 32.46522921557468)
    20 = @(38.87509062681747, 34.01866967722953)
    21 = @(39.70431054246355, 35.8059883368877)
    22 = @(44.91982769538773, 37.79632345265016)
    23 = @(48.29623832147909, 39.40986368251237)
    24 = @(51.01951624573073, 41.18438893606005)
    25 = @(51.37813468511153, 43.06230284459767)
    26 = @(52.74029464863034, 44.68744431235028)
    27 = @(54.50386636946551, 46.30555112148593)
    28 = @(54.87548333302324, 47.87172233768119)
    29 = @(57.44394973629721, 49.74663066837538)
    30 = @(58.19371234746801, 51.53826430886667)
    31 = @(61.05161607736118, 53.29391997002056)
    32 = @(66.58562711849606, 54.88388479987154)
    33 = @(74.80288205734669, 56.60695658848047)
    34 = @(75.74048677556678, 58.42355015082547)
    35 = @(77.64500024892949, 60.16795501091609)
    36 = @(78.56925677134853, 61.87360595942839)
    37 = @(78.81318145183432, 63.84052000165484)
    38 = @(78.99111476734312, 65.48939221638673)
    39 = @(79.75796682416843, 67.77190893764906)
    40 = @(82.17438814882016, 69.42334877475976)
    41 = @(82.59806165452208, 71.26117643576106)
    42 = @(84.11122107056853, 73.11046073806814)
    43 = @(86.94059514999302, 74.82242569141349)
    44 = @(89.40139508607268, 76.4434027924929)
    45 = @(90.82586510426854, 78.05868930261877)
    46 = @(92.46536369502256, 80.06676556254101)
    47 = @(94.2167820233657, 81.65322110015174)
    48 = @(94.68071828299796, 83.57905938231887)
    49 = @(98.52989125910187, 85.35676270936725)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item([int]$row, 2).Value = $pair[0]
    $ws.Cells.Item([int]$row, 3).Value = $pair[1]
}

# Row 50 (time = 48) no longer exists in the faster simulation output; remove it
# so the used range shrinks back to A1:C49.
$ws.Rows(50).Delete()
